$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the standalone "8452037 - Elisabeth..." row), shifting everything below up by one.
$ws.Rows.Item(13).Delete()

# Update cell contents that changed as a result of the restructuring.
$ws.Range("B10").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("C10").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("B14").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading'
$ws.Range("C14").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading'
$ws.Range("B15").Value = '01/01/2017'
$ws.Range("C15").Value = '01/01/2017'
$ws.Range("B16").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
$ws.Range("C16").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
$ws.Range("B18").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("C18").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("B19").Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Range("C19").Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Range("B20").Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Range("C20").Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Range("B23").Value = '8800009 -  Canto Coral I  (Requisito)
'
$ws.Range("C23").Value = '8800009 -  Canto Coral I  (Requisito)
'
